$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 252; $r++) {
    if ($r -le 29) {
        $ws.Cells.Item($r, 3).Value = 7318
    } else {
        $ws.Cells.Item($r, 3).Value = 7293
    }
}
